$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the short-url column (B) for all data rows (2-38) to the new value.
$ws.Range("B2:B38").Value = "M3tPQd"

# 2. Row 35: coo_id/coo_name/coo/coo_iso changes to India; swap asylum_seekers/returned_refugees.
$ws.Range("F35").Value = "'88"
$ws.Range("G35").Value = "India"
$ws.Range("H35").Value = "IND"
$ws.Range("I35").Value = "IND"
$ws.Range("N35").Value = "'0"
$ws.Range("O35").Value = "'5"

# 3. Row 37: coo_id/coo_name/coo/coo_iso changes to Tonga; swap asylum_seekers/returned_refugees.
$ws.Range("F37").Value = "'193"
$ws.Range("G37").Value = "Tonga"
$ws.Range("H37").Value = "TON"
$ws.Range("I37").Value = "TON"
$ws.Range("N37").Value = "'5"
$ws.Range("O37").Value = "'0"

# 4. Row 38: coo_id/coo_name/coo/coo_iso changes to Zimbabwe; swap asylum_seekers/returned_refugees.
$ws.Range("F38").Value = "'214"
$ws.Range("G38").Value = "Zimbabwe"
$ws.Range("H38").Value = "ZIM"
$ws.Range("I38").Value = "ZWE"
$ws.Range("N38").Value = "'5"
$ws.Range("O38").Value = "'0"

# 5. Delete row 39 (the Tonga/2024 duplicate row that was removed).
$ws.Rows("39").Delete()
